$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Change B10, B13, B23 from "Otros" to "Flow"
$ws.Range("B10").Value = "Flow"
$ws.Range("B13").Value = "Flow"
$ws.Range("B23").Value = "Flow"

# 2. Apply border + center alignment style to B10, then copy to B13/B23 (keeps the same style index)
$b10 = $ws.Range("B10")
$b10.HorizontalAlignment = -4108
$b10.Borders.Item(8).Weight = 2
$b10.Borders.Item(9).Weight = 2
$b10.Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("B23").PasteSpecial(-4122)

# 3. Apply AutoFilter for column B = "Flow"
$ws.Range("A1:F26").AutoFilter(2, @("Flow"), 7)

Write-Host "done"
